$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 3.900430680208489
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 654.126487776743
